$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 through 5 (extra data rows no longer needed)
$ws.Range("A3:C5").EntireRow.Delete()

# Update row 2 with new values
$ws.Range("A2").Value = "agnihotriaman@gmail.com"
$ws.Range("B2").Value = "124ef1"
$ws.Range("C2").Value = "testByCompanies"

# Add new column D with header (copy header formatting) and value
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "timeStamp"

$ws.Range("D2").Value = "2025-08-22 14:35:21"
